# The commit swaps the contents of ppt/theme/theme1.xml ("Integral" theme,
# used by the slide master / main deck) and ppt/theme/theme2.xml
# ("Office Theme", used by the notes master). The font scheme and format
# scheme blocks are byte-identical between the two theme parts - only the
# 12 colour-scheme entries (<a:clrScheme>) differ, so reproducing the swap
# is just a matter of re-pointing the deck's colour scheme to the colours
# that used to live in the "Office Theme" part.

function HexToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Colour order inside <a:clrScheme> == ThemeColorScheme.Item(1..12)
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$p = $ppt.ActivePresentation

$scheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $scheme.Count; $i++) {
    $scheme.Item($i).RGB = HexToRgbInt($officeThemeColors[$i - 1])
}
